# Borgs_Retrospective - Sprint 2.pptx update
# Fills in the placeholder numbered bullets ("1","2","3","4", "TEAM - ",
# "Carter -", ...) with the actual retrospective notes, on slides 1-3.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 1 - "What Went Well?"
# ---------------------------------------------------------------------
$s1  = $p.Slides.Item(1)
$sh1 = $s1.Shapes.Item(1)
$tr1 = $sh1.TextFrame.TextRange

# Paragraph 9: "Daniel -" -> "Daniel - We were able to grasp..."
# (text shares a prefix with the new text, so write a scratch value
# first to stop the host from keeping the shared prefix as its own run)
$para = $tr1.Paragraphs(9,1)
$para.Text = "X"
$para = $tr1.Paragraphs(9,1)
$para.Text = "Daniel – We were able to grasp a better understanding of how our server interacts with the clients."

# Paragraph 8: " Carter -" -> " Carter - Sprint was beneficial..."
$para = $tr1.Paragraphs(8,1)
$para.Text = "X"
$para = $tr1.Paragraphs(8,1)
$para.Text = " Carter – Sprint was beneficial because allowed for us to have a better understanding of our web application."

# Paragraph 7: " Walter -" is two runs (" Walte" + "r -"); only the
# second run gains a trailing space ("r -" -> "r - ").
$para = $tr1.Paragraphs(7,1)
$run2Start = $para.Start + 6
$run2Len   = $para.Length - 1 - 6
$run2 = $tr1.Characters($run2Start, $run2Len)
$run2.Text = "r – "

# Paragraphs 5 ("3") and 6 ("4") merge into a single paragraph with two
# runs: "Zenhub" + " was updated consistently".
$para6 = $tr1.Paragraphs(6,1)
$para6.Delete() | Out-Null
$para5 = $tr1.Paragraphs(5,1)
$para5.Text = "Zenhub"
$para5.InsertAfter(" was updated consistently") | Out-Null

# Paragraph 4: "2" -> "Everyone had a better understanding of the project"
$para = $tr1.Paragraphs(4,1)
$para.Text = "Everyone had a better understanding of the project"

# Paragraph 3: "1" -> "We planned accordingly." and a brand-new bullet
# paragraph is added right after it: "Broke down tasks into more
# manageable tasks." (same paragraph formatting as paragraph 3).
$para3 = $tr1.Paragraphs(3,1)
$para3.InsertAfter("`r") | Out-Null
$para3 = $tr1.Paragraphs(3,1)
$para3.Text = "We planned accordingly."
$para3new = $tr1.Paragraphs(4,1)
$para3new.Text = "Broke down tasks into more manageable tasks."

# ---------------------------------------------------------------------
# Slide 2 - "What Might Be Impeding Us from Performing Better?"
# ---------------------------------------------------------------------
$s2  = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(1)
$tr2 = $sh2.TextFrame.TextRange

# Paragraph 4: "3" run is removed entirely; the paragraph becomes empty.
$para = $tr2.Paragraphs(4,1)
$para.Text = ""

# Paragraph 3: "2" -> "Difficulty meeting up with group due to each other's schedules"
$para = $tr2.Paragraphs(3,1)
$para.Text = "Difficulty meeting up with group due to each other's schedules"

# Paragraph 2: "1" -> "Lack of time due to other responsibilities. "
$para = $tr2.Paragraphs(2,1)
$para.Text = "Lack of time due to other responsibilities. "

# ---------------------------------------------------------------------
# Slide 3 - "What Can We do to Improve?"
# ---------------------------------------------------------------------
$s3  = $p.Slides.Item(3)
$sh3 = $s3.Shapes.Item(1)
$tr3 = $sh3.TextFrame.TextRange

# Paragraph 4: "Carter -" -> "Carter - More group meetups/group chat..."
$para = $tr3.Paragraphs(4,1)
$para.Text = "X"
$para = $tr3.Paragraphs(4,1)
$para.Text = "Carter – More group meetups/group chat sessions to discuss issues."

# Paragraph 2: "TEAM - " -> "TEAM - More communication between the whole team."
$para = $tr3.Paragraphs(2,1)
$para.Text = "X"
$para = $tr3.Paragraphs(2,1)
$para.Text = "TEAM – More communication between the whole team."
